{"js": "// Replace each three-digit division expression with its updated value.\n// Each old expression occurs exactly once in the document, so a targeted\n// search-and-replace for every pair reproduces the diff exactly.\nconst replacements = [\n  [\"493\u00f72=\", \"276\u00f77=\"],\n  [\"799\u00f79=\", \"854\u00f78=\"],\n  [\"266\u00f73=\", \"811\u00f74=\"],\n  [\"185\u00f78=\", \"346\u00f77=\"],\n  [\"663\u00f74=\", \"963\u00f75=\"],\n  [\"720\u00f79=\", \"931\u00f78=\"],\n  [\"769\u00f72=\", \"143\u00f78=\"],\n  [\"198\u00f74=\", \"715\u00f75=\"],\n  [\"553\u00f78=\", \"234\u00f78=\"],\n  [\"541\u00f72=\", \"563\u00f77=\"],\n  [\"236\u00f77=\", \"188\u00f78=\"],\n  [\"817\u00f78=\", \"278\u00f74=\"],\n  [\"397\u00f76=\", \"724\u00f75=\"],\n  [\"873\u00f75=\", \"757\u00f74=\"],\n  [\"169\u00f77=\", \"609\u00f75=\"],\n  [\"615\u00f77=\", \"219\u00f72=\"],\n  [\"301\u00f74=\", \"781\u00f76=\"],\n  [\"670\u00f74=\", \"259\u00f73=\"],\n  [\"856\u00f72=\", \"321\u00f76=\"],\n  [\"102\u00f72=\", \"360\u00f77=\"],\n  [\"448\u00f76=\", \"690\u00f73=\"],\n  [\"745\u00f76=\", \"812\u00f79=\"],\n  [\"355\u00f78=\", \"333\u00f72=\"],\n  [\"801\u00f73=\", \"223\u00f72=\"],\n  [\"603\u00f74=\", \"841\u00f72=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each three-digit division expression in the practice table\n# with its new value. Every old expression occurs exactly once in the\n# document, so a Find/Replace pass per pair reproduces the diff exactly.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"493\u00f72=\", \"276\u00f77=\"),\n    @(\"799\u00f79=\", \"854\u00f78=\"),\n    @(\"266\u00f73=\", \"811\u00f74=\"),\n    @(\"185\u00f78=\", \"346\u00f77=\"),\n    @(\"663\u00f74=\", \"963\u00f75=\"),\n    @(\"720\u00f79=\", \"931\u00f78=\"),\n    @(\"769\u00f72=\", \"143\u00f78=\"),\n    @(\"198\u00f74=\", \"715\u00f75=\"),\n    @(\"553\u00f78=\", \"234\u00f78=\"),\n    @(\"541\u00f72=\", \"563\u00f77=\"),\n    @(\"236\u00f77=\", \"188\u00f78=\"),\n    @(\"817\u00f78=\", \"278\u00f74=\"),\n    @(\"397\u00f76=\", \"724\u00f75=\"),\n    @(\"873\u00f75=\", \"757\u00f74=\"),\n    @(\"169\u00f77=\", \"609\u00f75=\"),\n    @(\"615\u00f77=\", \"219\u00f72=\"),\n    @(\"301\u00f74=\", \"781\u00f76=\"),\n    @(\"670\u00f74=\", \"259\u00f73=\"),\n    @(\"856\u00f72=\", \"321\u00f76=\"),\n    @(\"102\u00f72=\", \"360\u00f77=\"),\n    @(\"448\u00f76=\", \"690\u00f73=\"),\n    @(\"745\u00f76=\", \"812\u00f79=\"),\n    @(\"355\u00f78=\", \"333\u00f72=\"),\n    @(\"801\u00f73=\", \"223\u00f72=\"),\n    @(\"603\u00f74=\", \"841\u00f72=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $find.Text,       # FindText\n        $false,           # MatchCase\n        $false,           # MatchWholeWord\n        $false,           # MatchWildcards\n        $false,           # MatchSoundsLike\n        $false,           # MatchAllWordForms\n        $true,            # Forward\n        1,                # Wrap (wdFindContinue)\n        $false,           # Format\n        $find.Replacement.Text,  # ReplaceWith\n        2                 # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
